# Updates cryptos list price (D) and 1h volume-change (E) columns per the
# 2024-09-20 04:57:29 UTC GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimals (e.g. "1.00", "0.0532") that Excel
# would otherwise auto-convert to numbers; a leading apostrophe forces them
# to stay text, matching the column's existing string formatting.

$ws.Range("D2").Value = '63.855.96'
$ws.Range("E2").Value = '  +2.94%  '
$ws.Range("D3").Value = '2.539.59'
$ws.Range("E3").Value = '  +5.42%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''573.88'  # was '573.86', force text
$ws.Range("E5").Value = '  +2.17%  '
$ws.Range("D6").Value = '''148.40'  # was '147.79', force text
$ws.Range("E6").Value = '  +6.79%  '
$ws.Range("E8").Value = '  +0.63%  '
$ws.Range("D9").Value = '2.539.57'
$ws.Range("E9").Value = '  +5.51%  '
$ws.Range("E10").Value = '  +2.51%  '
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("E12").Value = '  +1.69%  '
$ws.Range("E13").Value = '  +3.34%  '
$ws.Range("D14").Value = '''28.04'  # was '27.87', force text
$ws.Range("E14").Value = '  +8.98%  '
$ws.Range("D15").Value = '2.994.22'
$ws.Range("E15").Value = '  +5.61%  '
$ws.Range("D16").Value = '63.678.16'
$ws.Range("E16").Value = '  +2.79%  '
$ws.Range("E17").Value = '  +3.85%  '
$ws.Range("D18").Value = '2.539.19'
$ws.Range("E18").Value = '  +5.62%  '
$ws.Range("D19").Value = '''11.57'  # was '11.52', force text
$ws.Range("E19").Value = '  +5.04%  '
$ws.Range("D20").Value = '''344.10'  # was '343.93', force text
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("D21").Value = '''4.37'  # was '4.36', force text
$ws.Range("E21").Value = '  +3.52%  '
$ws.Range("D22").Value = '''6.90'  # was '6.89', force text
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("D24").Value = '''65.95'  # was '66.12', force text
$ws.Range("E24").Value = '  +1.72%  '
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("E26").Value = '  +5.16%  '
$ws.Range("D27").Value = '''1.00'  # was '0.999', force text
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '''8.29'  # was '8.26', force text
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").Value = '''1.43'  # was '1.42', force text
$ws.Range("E29").Value = '  +3.44%  '
$ws.Range("E30").Value = '  +7.27%  '
$ws.Range("E31").Value = '  +4.04%  '
$ws.Range("E32").Value = '  +7.43%  '
$ws.Range("D33").Value = '''177.05'  # was '177.06', force text
$ws.Range("E33").Value = '  +3.68%  '
$ws.Range("D34").Value = '''1.57'  # was '1.55', force text
$ws.Range("E34").Value = '  +12.04%  '
$ws.Range("D35").Value = '''423.46'  # was '416.17', force text
$ws.Range("E35").Value = '  +16.77%  '
$ws.Range("D36").Value = '''0.407'  # was '0.405', force text
$ws.Range("E36").Value = '  +3.59%  '
$ws.Range("D37").Value = '''19.17'  # was '19.13', force text
$ws.Range("E37").Value = '  +3.74%  '
$ws.Range("D38").Value = '''4.47'  # was '4.45', force text
$ws.Range("E38").Value = '  -2.19%  '
$ws.Range("E40").Value = '  +6.07%  '
$ws.Range("D41").Value = '''1.00'  # was '0.999', force text
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").Value = '''40.75'  # was '40.78', force text
$ws.Range("E42").Value = '  +4.53%  '
$ws.Range("E43").Value = '  +6.33%  '
$ws.Range("D45").Value = '''20.96'  # was '20.95', force text
$ws.Range("E45").Value = '  +2.61%  '
$ws.Range("E46").Value = '  +4.40%  '
$ws.Range("D47").Value = '''0.0532'  # was '0.0531', force text
$ws.Range("E47").Value = '  +2.68%  '
$ws.Range("D48").Value = '''0.0968'  # was '0.0966', force text
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("D49").Value = '''19.01'  # was '19.04', force text
$ws.Range("E49").Value = '  +6.61%  '
$ws.Range("E50").Value = '  +5.19%  '
$ws.Range("D51").Value = '''1.83'  # was '1.84', force text
$ws.Range("E51").Value = '  +9.60%  '
